$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row updates: (row, country, B, C, D, E, F, G, H)
$updates = @(
    ,@(4, 'Estados Unidos', 1139263, 8233, 162100, 910757, 16481, 653, 66406)
    ,@(10, 'Turquia', 124375, 1983, 58259, 62780, 1445, 78, 3336)
    ,@(11, 'Rusia', 124054, 9623, 15013, 107819, 2300, 53, 1222)
    ,@(21, 'Ecuador', 27464, 1128, 2132, 23961, 149, 308, 1371)
    ,@(28, 'Chile', 18435, 1427, 9572, 8616, 425, 13, 247)
    ,@(51, 'Egipto', 6193, 298, 1522, 4256, 0, 9, 415)
    ,@(52, 'Malasia', 6176, 105, 4326, 1747, 31, 0, 103)
    ,@(53, 'Sudafrica', 5951, 0, 2382, 3453, 36, 0, 116)
    ,@(55, 'Marruecos', 4729, 160, 1256, 3300, 1, 2, 173)
    ,@(104, 'Somalia', 671, 70, 517, 606, 2, 3, 31)
    ,@(105, 'Burkina Faso', 649, 0, 517, 88, 0, 0, 44)
    ,@(106, 'Uruguay', 648, 0, 435, 196, 10, 0, 17)
    ,@(107, 'Guatemala', 644, 45, 72, 556, 5, 0, 16)
    ,@(108, 'Consejo Danes para los Refugiados', 604, 0, 75, 497, 2, 0, 32)
    ,@(193, 'Namibia', 16, 0, 8, 8, 0, 0, 0)
    ,@(194, 'San Vicente y las Granadinas', 16, 0, 8, 8, 0, 0, 0)
    ,@(217, 'Comoras', 1, 0, 0, 1, 0, 0, 0)
    ,@(218, 'San Pedro y Miquelon', 1, 0, 0, 1, 0, 0, 0)
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 1).Value = $u[1]
    for ($c = 2; $c -le 8; $c++) {
        $ws.Cells.Item($row, $c).Value = $u[$c]
    }
}

Write-Host "Update complete."